$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.974.15'
$ws.Range("E2").Value = '  -1.98%  '
$ws.Range("D3").Value = '2.972.33'
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("D5").Value = '''594.13'
$ws.Range("E5").Value = '  +1.71%  '
$ws.Range("D6").Value = '''142.23'
$ws.Range("E6").Value = '  -2.05%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -1.39%  '
$ws.Range("D9").Value = '2.969.65'
$ws.Range("E9").Value = '  -0.57%  '
$ws.Range("D10").Value = '''0.146'
$ws.Range("E10").Value = '  -1.12%  '
$ws.Range("D11").Value = '''5.95'
$ws.Range("E11").Value = '  +4.30%  '
$ws.Range("D12").Value = '''0.450'
$ws.Range("E12").Value = '  +2.19%  '
$ws.Range("E13").Value = '  -0.68%  '
$ws.Range("D14").Value = '''33.96'
$ws.Range("E14").Value = '  -1.38%  '
$ws.Range("E15").Value = '  +2.34%  '
$ws.Range("D16").Value = '3.459.01'
$ws.Range("E16").Value = '  -0.64%  '
$ws.Range("D17").Value = '61.073.09'
$ws.Range("E17").Value = '  -1.85%  '
$ws.Range("D18").Value = '''6.83'
$ws.Range("E18").Value = '  -1.99%  '
$ws.Range("D19").Value = '2.966.10'
$ws.Range("E19").Value = '  -0.82%  '
$ws.Range("D20").Value = '''449.02'
$ws.Range("E20").Value = '  -1.46%  '
$ws.Range("D21").Value = '''14.11'
$ws.Range("E21").Value = '  +2.34%  '
$ws.Range("D22").Value = '''0.677'
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("D23").Value = '''7.24'
$ws.Range("E23").Value = '  -0.98%  '
$ws.Range("D24").Value = '''82.03'
$ws.Range("E24").Value = '  +2.66%  '
$ws.Range("E25").Value = '  -4.31%  '
$ws.Range("D26").Value = '''10.31'
$ws.Range("E26").Value = '  +2.55%  '
$ws.Range("D27").Value = '''11.86'
$ws.Range("E27").Value = '  -2.45%  '
$ws.Range("E28").Value = '  +0.22%  '
$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("E30").Value = '  +1.87%  '
$ws.Range("D31").Value = '''7.07'
$ws.Range("E31").Value = '  -0.42%  '
$ws.Range("E32").Value = '  -1.85%  '
$ws.Range("D33").Value = '''27.21'
$ws.Range("E33").Value = '  +2.02%  '
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("D35").Value = '0.0₃0802'
$ws.Range("E35").Value = '  +3.02%  '
$ws.Range("E36").Value = '  -1.17%  '
$ws.Range("D37").Value = '''5.75'
$ws.Range("E37").Value = '  +0.88%  '
$ws.Range("D38").Value = '''49.98'
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '''2.02'
$ws.Range("E39").Value = '  -3.19%  '
$ws.Range("B40").Value = 'Cosmos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D40").Value = '''8.94'
$ws.Range("E40").Value = '  +0.61%  '
$ws.Range("E41").Value = '  +6.50%  '
$ws.Range("E42").Value = '  -3.62%  '
$ws.Range("D43").Value = '''385.15'
$ws.Range("E43").Value = '  -4.47%  '
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").Value = '''0.266'
$ws.Range("E44").Value = '  -2.59%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = '''0.0346'
$ws.Range("E45").Value = '  -0.67%  '
$ws.Range("B46").Value = 'Arweave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D46").Value = '''38.42'
$ws.Range("E46").Value = '  -1.67%  '
$ws.Range("D47").Value = '2.689.66'
$ws.Range("E47").Value = '  -2.27%  '
$ws.Range("D48").Value = '''129.48'
$ws.Range("E48").Value = '  +1.29%  '
$ws.Range("E50").Value = '  -0.64%  '
$ws.Range("D51").Value = '''2.12'
$ws.Range("E51").Value = '  -0.22%  '
